# Update harmonized emissions values on the "data" sheet to reflect the
# recalculated harmonization results (base/harmonize year change), per
# the regression test baseline update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: prefix|Emissions|BC|Harmonized-DB
$ws.Range("G2").Value = 32.33588516746411
$ws.Range("H2").Value = 34.01818181818182
$ws.Range("I2").Value = 35.52057416267942
$ws.Range("J2").Value = 36.84306220095694
$ws.Range("K2").Value = 37.98564593301435
$ws.Range("L2").Value = 38.94832535885168
$ws.Range("M2").Value = 39.7311004784689
$ws.Range("N2").Value = 40.33397129186603
$ws.Range("O2").Value = 40.75693779904307

# Row 3: prefix|Emissions|BC|sector1|Harmonized-DB
$ws.Range("G3").Value = 15.16842105263158
$ws.Range("H3").Value = 16.04210526315789
$ws.Range("I3").Value = 16.83157894736842
$ws.Range("J3").Value = 17.53684210526316
$ws.Range("K3").Value = 18.15789473684211
$ws.Range("L3").Value = 18.69473684210526
$ws.Range("M3").Value = 19.14736842105263
$ws.Range("N3").Value = 19.51578947368421
$ws.Range("O3").Value = 19.8

# Row 4: prefix|Emissions|BC|sector2|Harmonized-DB
$ws.Range("G4").Value = 17.16746411483254
$ws.Range("H4").Value = 17.97607655502392
$ws.Range("I4").Value = 18.68899521531101
$ws.Range("J4").Value = 19.30622009569378
$ws.Range("K4").Value = 19.82775119617225
$ws.Range("L4").Value = 20.25358851674642
$ws.Range("M4").Value = 20.58373205741627
$ws.Range("N4").Value = 20.81818181818182
$ws.Range("O4").Value = 20.95693779904306
